# Apply updated cryptocurrency price/volume data to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimal numbers (e.g. "706.19"). Setting
# .Value directly on those cells would let Excel auto-convert the text into
# a number, which would change the stored cell type away from the original
# inline string. Force those specific cells to Text format first, write the
# string values, then restore each cell's original style so no stray
# formatting is left behind.
$numericPriceCells = @("D5", "D6", "D11", "D13", "D14", "D18", "D19", "D21", "D22", "D23", "D24", "D28", "D33", "D35", "D36", "D38", "D40", "D42", "D43", "D46", "D47", "D48", "D49", "D50")
$savedStyles = @{}
foreach ($addr in $numericPriceCells) {
    $r = $ws.Range($addr)
    $savedStyles[$addr] = $r.Style
    $r.NumberFormat = "@"
}

$ws.Range('D2').Value = '71.097.90'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '3.811.87'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '706.19'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '171.38'
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('D7').Value = '3.810.50'
$ws.Range('E7').Value = '  -1.26%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('D11').Value = '7.74'
$ws.Range('E11').Value = '  +7.10%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = '0.0000252'
$ws.Range('E13').Value = '  -2.10%  '
$ws.Range('D14').Value = '35.87'
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('D15').Value = '4.454.85'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '3.806.59'
$ws.Range('E16').Value = '  -4.15%  '
$ws.Range('D17').Value = '71.071.05'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '7.15'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '17.43'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '502.27'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = '10.73'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').Value = '0.724'
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('D24').Value = '84.21'
$ws.Range('E24').Value = '  -0.98%  '
$ws.Range('E25').Value = '  -3.42%  '
$ws.Range('D26').Value = '3.963.18'
$ws.Range('E26').Value = '  -1.19%  '
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').Value = '10.38'
$ws.Range('E28').Value = '  -2.57%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('E30').Value = '  -4.26%  '
$ws.Range('E31').Value = '  -3.86%  '
$ws.Range('E32').Value = '  -0.57%  '
$ws.Range('D33').Value = '7.34'
$ws.Range('E33').Value = '  -2.74%  '
$ws.Range('E34').Value = '  -1.89%  '
$ws.Range('D35').Value = '0.173'
$ws.Range('E35').Value = '  -4.52%  '
$ws.Range('D36').Value = '9.16'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('D37').Value = '3.777.66'
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('D38').Value = '0.998'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').Value = '  -2.60%  '
$ws.Range('D40').Value = '2.37'
$ws.Range('E40').Value = '  -1.00%  '
$ws.Range('E41').Value = '  -2.92%  '
$ws.Range('D42').Value = '5.93'
$ws.Range('E42').Value = '  -1.79%  '
$ws.Range('D43').Value = '3.26'
$ws.Range('E43').Value = '  -4.53%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').Value = '166.96'
$ws.Range('E46').Value = '  +1.96%  '
$ws.Range('D47').Value = '0.000313'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').Value = '49.06'
$ws.Range('E48').Value = '  -0.33%  '
$ws.Range('D49').Value = '422.36'
$ws.Range('E49').Value = '  +1.40%  '
$ws.Range('D50').Value = '8.61'
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('E51').Value = '  -2.21%  '

# Restore original styles on the cells we temporarily switched to Text format.
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).Style = $savedStyles[$addr]
}
